# Update countries & provincias Spain
# Applies the May-27-2020 01:05 data refresh to the "Pais" worksheet:
#  - in-place numeric updates for several countries whose rank didn't change
#  - Colombia's case counts grew enough to overtake Kuwait & Polonia (rows 36-38)
#  - Nigeria's case counts grew enough to overtake Oman (rows 58-59)
#  - Groenlandia and Islas Turcas y Caicos swap places (rows 207-208)
#  - footer timestamp text updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 01:05"

# --- Row 4: Estados Unidos (rank 8, unchanged position) ---
$ws.Range("B4").Value = 1724504
$ws.Range("C4").Value = 18278
$ws.Range("D4").Value = 477504
$ws.Range("E4").Value = 1146490
$ws.Range("G4").Value = 705
$ws.Range("H4").Value = 100510

# --- Row 5: Brasil (rank 9, unchanged position) ---
$ws.Range("B5").Value = 391222
$ws.Range("C5").Value = 14553
$ws.Range("D5").Value = 158593
$ws.Range("E5").Value = 208117
$ws.Range("G5").Value = 990
$ws.Range("H5").Value = 24512

# --- Row 15: Peru (rank 19, unchanged position) ---
$ws.Range("B15").Value = 129751
$ws.Range("C15").Value = 5772
$ws.Range("D15").Value = 52906
$ws.Range("E15").Value = 73057
$ws.Range("G15").Value = 159
$ws.Range("H15").Value = 3788

# --- Rows 36-38: Colombia overtakes Kuwait and Polonia ---
# Row 36 becomes Colombia with its updated (increased) counts.
$ws.Range("A36").Value = "Colombia"
$ws.Range("B36").Value = 23003
$ws.Range("C36").Value = 1022
$ws.Range("D36").Value = 5511
$ws.Range("E36").Value = 16716
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 26
$ws.Range("H36").Value = 776

# Row 37 becomes Kuwait (its data is unchanged, just shifted down one row).
$ws.Range("A37").Value = "Kuwait"
$ws.Range("B37").Value = 22575
$ws.Range("C37").Value = 608
$ws.Range("D37").Value = 7306
$ws.Range("E37").Value = 15097
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 7
$ws.Range("H37").Value = 172

# Row 38 becomes Polonia (its data is unchanged, just shifted down one row).
$ws.Range("A38").Value = "Polonia"
$ws.Range("B38").Value = 22074
$ws.Range("C38").Value = 443
$ws.Range("D38").Value = 10020
$ws.Range("E38").Value = 11030
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 17
$ws.Range("H38").Value = 1024

# --- Row 54: Chequia (rank 58, unchanged position) ---
$ws.Range("B54").Value = 9050
$ws.Range("C54").Value = 48
$ws.Range("D54").Value = 6270
$ws.Range("E54").Value = 2463

# --- Row 57: Noruega (rank 61, unchanged position) ---
$ws.Range("B57").Value = 8383
$ws.Range("C57").Value = 19
$ws.Range("E57").Value = 421

# --- Rows 58-59: Nigeria overtakes Oman ---
# Row 58 becomes Nigeria with its updated (increased) counts.
$ws.Range("A58").Value = "Nigeria"
$ws.Range("B58").Value = 8344
$ws.Range("C58").Value = 276
$ws.Range("D58").Value = 2385
$ws.Range("E58").Value = 5710
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 16
$ws.Range("H58").Value = 249

# Row 59 becomes Oman (its data is unchanged, just shifted down one row).
$ws.Range("A59").Value = "Oman"
$ws.Range("B59").Value = 8118
$ws.Range("C59").Value = 348
$ws.Range("D59").Value = 2067
$ws.Range("E59").Value = 6014
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 37

# --- Rows 207-208: Groenlandia and Islas Turcas y Caicos swap places ---
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("D207").Value = 11
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 10
$ws.Range("H208").Value = 1
